$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of this data block (row 232),
# pushing the existing rows 232:298 down to 233:299.
$ws.Rows(232).Insert()

# Populate the newly inserted row 232 with the new record's data.
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44855
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100108
$ws.Range("H232").Value = "Tropicales y subtropicales"
$ws.Range("I232").Value = 100108005
$ws.Range("J232").Value = "Piña"
$ws.Range("K232").Value = "Caramelo"
$ws.Range("L232").Value = "Tercera"
$ws.Range("M232").Value = 280
$ws.Range("N232").Value = 19000
$ws.Range("O232").Value = 19000
$ws.Range("P232").Value = 19000
$ws.Range("Q232").Value = "$/caja 16 unidades"
$ws.Range("R232").Value = "Ecuador"
$ws.Range("S232").Value = 1188
$ws.Range("T232").Value = 16
